# Apply the changes described in the commit "updated 4.0 files and mdl"
# to the "Maximum Capacity Factor.xlsx" workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet: bump the "last updated" date in C1
#    (45320 -> 45392, serial date 1/29/2024 -> 4/10/2024)
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# ---------------------------------------------------------------------
# 2. "MCF" sheet: every Maximum Capacity Factor that was 0.85 or 0.95
#    is raised to 1 (full capacity factor). Cells already at 0 or 1
#    are left untouched.
# ---------------------------------------------------------------------
$wsMcf = $wb.Worksheets.Item("MCF")

$wsMcf.Range("B2").Value = 1
$wsMcf.Range("B3").Value = 1
$wsMcf.Range("B4").Value = 1
$wsMcf.Range("B6").Value = 1
$wsMcf.Range("B10").Value = 1
$wsMcf.Range("B11").Value = 1
$wsMcf.Range("B12").Value = 1
$wsMcf.Range("B13").Value = 1
$wsMcf.Range("B14").Value = 1
$wsMcf.Range("B16").Value = 1
$wsMcf.Range("B17").Value = 1
$wsMcf.Range("B18").Value = 1

# B19, B20, B21, B22, B24, B25 hold formulas referencing the cells
# above (=B2, =B4, =B10, =B14, =B4, =B4); they recalculate automatically
# once their source cells change, but force a recalculation to be sure
# the cached <v> values in the saved file reflect the new results.
$excel.CalculateFullRebuild()
$wb.Application.Calculate()

# ---------------------------------------------------------------------
# 3. Update the selection left behind on the MCF sheet (E8 -> B17)
# ---------------------------------------------------------------------
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
